$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coefficient values for rows 2-4
$ws.Range("B2").Value = -0.4810784687029628
$ws.Range("B3").Value = 0.9476221761610764
$ws.Range("B4").Value = 225.2805515915206

# Remove row 5 entirely (the "4" / 1.5720823297346278 row)
$ws.Range("A5:B5").EntireRow.Delete()
